$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate data rows 2-6 up by one, with old row 2 wrapping around to row 6 ---
# (old row3->new row2, old row4->new row3, old row5->new row4, old row6->new row5,
#  old row2->new row6)

# 1) Stash old row 2 (values + styles) far away so later row-shifts don't disturb it.
$ws.Range("A2:B2").Copy($ws.Range("A100:B100"))

# 2) Delete row 2 entirely; rows 3-6 shift up to 2-5 automatically (values + styles
#    come along for the ride).
$ws.Rows(2).Delete()

# 3) The stashed data (originally placed at row 100) has itself shifted up to row 99
#    because of the deletion above. Copy it down into the now-empty row 6.
$ws.Range("A99:B99").Copy($ws.Range("A6:B6"))

# 4) Clean up the temporary holding cells.
$ws.Range("A99:B99").Clear()

# At this point rows 2-6 already carry the correct values AND the correct cell
# styles (the bordered "hyperlink look" followed the row shift). Remember which
# style each of A2/A3/A6 currently has, and A3's real text, so both can be
# restored after the hyperlink rebuild below (adding a hyperlink can otherwise
# spawn near-duplicate styles, and a TextToDisplay overwrites the cell text).
$ws.Range("A2").Copy($ws.Range("C100"))
$ws.Range("A3").Copy($ws.Range("C101"))
$ws.Range("A6").Copy($ws.Range("C102"))
$a3Text = $ws.Range("A3").Value2

# --- Hyperlinks are not re-targeted automatically by the row delete/shift, so
#     rebuild the three mailto hyperlinks from scratch at their final cells. ---
$ws.Range("A1:B200").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:sasikala.ars@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:sasikala.ars@gmail.com", "", "", "sasikala.ars@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:sasikala.ars@gmail.com")

# The TextToDisplay passed above overwrote A3's cell text; put the real text
# back (this does not disturb the "display" attribute already recorded on the
# hyperlink relationship).
$ws.Range("A3").Value = $a3Text

# Reassert the exact original styles on A2/A3/A6 (values are untouched by this;
# xlPasteFormats only copies formatting, not cell content).
$ws.Range("C100").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("C101").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("C102").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Remove the scratch cells used above.
$ws.Range("C100:C102").Clear()

# --- Update the selected cell shown when the sheet is opened ---
$ws.Range("A3").Select()
